$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row of data (row 6) mirroring the structure of the existing rows.
$row = 6

$ws.Cells.Item($row, 1).Value = 42588.471562500003
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = "Noun"

$ws.Cells.Item($row, 3).Value = 10858
$ws.Cells.Item($row, 4).Value = 6710
$ws.Cells.Item($row, 5).Value = 1264
$ws.Cells.Item($row, 6).Value = 121
$ws.Cells.Item($row, 7).Value = 67
$ws.Cells.Item($row, 8).Value = 64
$ws.Cells.Item($row, 9).Value = 35
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 2
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 100
